$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.566.00"
$ws.Range("E2").Value = "  -6.31%  "
$ws.Range("D3").Value = "3.250.87"
$ws.Range("E3").Value = "  -9.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "175.58"
$ws.Range("E5").Value = "  -12.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "505.29"
$ws.Range("E6").Value = "  -10.77%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.334.97"
$ws.Range("E7").Value = "  -6.52%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.580"
$ws.Range("E8").Value = "  -5.04%  "
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("E10").Value = "  -10.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.32"
$ws.Range("E11").Value = "  -5.76%  "
$ws.Range("E12").Value = "  -12.73%  "
$ws.Range("E13").Value = "  -11.02%  "
$ws.Range("E14").Value = "  -13.02%  "
$ws.Range("D15").Value = "3.743.49"
$ws.Range("E15").Value = "  -9.80%  "
$ws.Range("E16").Value = "  -4.80%  "
$ws.Range("D17").Value = "3.233.76"
$ws.Range("E17").Value = "  -9.61%  "
$ws.Range("D18").Value = "63.175.50"
$ws.Range("E18").Value = "  -6.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.78"
$ws.Range("E19").Value = "  -11.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.59"
$ws.Range("E20").Value = "  -12.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.925"
$ws.Range("E21").Value = "  -12.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "363.48"
$ws.Range("E22").Value = "  -9.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "78.73"
$ws.Range("E23").Value = "  -6.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.54"
$ws.Range("E24").Value = "  -14.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.56"
$ws.Range("E25").Value = "  -15.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.97"
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.69"
$ws.Range("E27").Value = "  -4.76%  "
$ws.Range("E28").Value = "  -10.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.90"
$ws.Range("E29").Value = "  -11.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.14"
$ws.Range("E30").Value = "  -11.08%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.80"
$ws.Range("E31").Value = "  -11.33%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "630.01"
$ws.Range("E32").Value = "  -5.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.47"
$ws.Range("E33").Value = "  -15.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.81"
$ws.Range("E34").Value = "  -10.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.17"
$ws.Range("E35").Value = "  -8.01%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.100"
$ws.Range("E37").Value = "  -10.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "34.79"
$ws.Range("E38").Value = "  -14.97%  "
$ws.Range("E39").Value = "  -10.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.993"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("E41").Value = "  -10.45%  "
$ws.Range("D42").Value = "2.770.65"
$ws.Range("E42").Value = "  -13.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.63"
$ws.Range("E43").Value = "  -18.04%  "
$ws.Range("E44").Value = "  -18.46%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.67"
$ws.Range("E45").Value = "  +17.46%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.55"
$ws.Range("E46").Value = "  -8.84%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0372"
$ws.Range("E47").Value = "  -8.91%  "
$ws.Range("E48").Value = "  -6.89%  "
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.21"
$ws.Range("E49").Value = "  -17.42%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.64"
$ws.Range("E50").Value = "  -5.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.75"
$ws.Range("E51").Value = "  -10.48%  "
